# Auto-generated edit script for fixtures.xlsx - "Common: Inventory looks good"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translations - Common")

# 1) Extend formatting (wrap-text style used by column A/B/C in this sheet) down into the new rows
#    by copying the format of the last existing data row (735) into the new block (736:761).
$ws.Range("A735:C735").Copy()
$ws.Range("A736:C761").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 2) Seed the new unique strings (in scratch cells, far away from the used range) in the exact
#    order they should land in the shared-string table, so the final workbook reproduces the
#    same append order as the authored edit. The scratch cells are cleared afterwards; any
#    strings that end up used by the real target cells below survive, others are pruned on save.
$ws.Range("Z1").Value = "inventory.aroma.index.title"
$ws.Range("Z2").Value = "Vaše aromata"
$ws.Range("Z3").Value = "inventory.base.index.title"
$ws.Range("Z4").Value = "Vaše báze"
$ws.Range("Z5").Value = "inventory.booster.index.title"
$ws.Range("Z6").Value = "Vaše boostery"
$ws.Range("Z7").Value = "inventory.mod.index.title"
$ws.Range("Z8").Value = "Vaše mody"
$ws.Range("Z9").Value = "inventory.cell.index.title"
$ws.Range("Z10").Value = "Vaše články"
$ws.Range("Z11").Value = "lab.cell.inventory.delete.success"
$ws.Range("Z12").Value = "Vybrané články byly úspěšně odstraněny."
$ws.Range("Z13").Value = "market.inventory.menu"
$ws.Range("Z14").Value = "lab.atomizer.inventory.delete.modal.title"
$ws.Range("Z15").Value = "Odstranit vybrané atomizéry?"
$ws.Range("Z16").Value = "lab.atomizer.inventory.delete.modal.content"
$ws.Range("Z17").Value = "<p>`n`tOpravdu si přejte odstranit vybrané atomizéry?`n</p>`n<p>`n`tNutno podotknout, že tato akce je <strong>silně destruktivní</strong> a není možné ji vzít zpět.`n</p>`n<p>`n`t<strong>Veškerá přidružená data budou smazána, včetně údajů o vapování, komentáře a další!</strong>`n</p>"
$ws.Range("Z18").Value = "lab.atomizer.inventory.delete.success"
$ws.Range("Z19").Value = "Vybrané atomizéry byly odstraněny. Bůh s vámi."
$ws.Range("Z20").Value = "lab.mod.inventory.delete.modal.title"
$ws.Range("Z21").Value = "lab.mod.inventory.delete.modal.content"
$ws.Range("Z22").Value = "<p>`n`tOpravdu si přejte odstranit vybrané mody?`n</p>`n<p>`n`tNutno podotknout, že tato akce je <strong>silně destruktivní</strong> a není možné ji vzít zpět.`n</p>`n<p>`n`t<strong>Veškerá přidružená data budou smazána, včetně údajů o vapování, komentáře a další!</strong>`n</p>"
$ws.Range("Z23").Value = "Vybrané mody byly odstraněny."
$ws.Range("Z24").Value = "lab.cell.inventory.delete.modal.title"
$ws.Range("Z25").Value = "Odstranit vybrané mody?"
$ws.Range("Z26").Value = "Odstranit vybrané články?"
$ws.Range("Z27").Value = "lab.cell.inventory.delete.modal.content"
$ws.Range("Z28").Value = "<p>`n`tOpravdu si přejte odstranit vybrané články?`n</p>`n<p>`n`tNutno podotknout, že tato akce je <strong>silně destruktivní</strong> a není možné ji vzít zpět.`n</p>`n<p>`n`t<strong>Veškerá přidružená data budou smazána, včetně údajů o vapování, komentáře a další!</strong>`n</p>"
$ws.Range("Z29").Value = "lab.mod.inventory.delete.success"
$ws.Range("Z30").Value = "Vybrané články byly odstraněny."
$ws.Range("Z31").Value = "inventory.build.menu"
$ws.Range("Z32").Value = "inventory.wire.menu"
$ws.Range("Z33").Value = "inventory.wire.index.title"
$ws.Range("Z34").Value = "Vaše odporové dráty"
$ws.Range("Z35").Value = "lab.wire.list.empty.title"
$ws.Range("Z36").Value = "Nemáte zakoupené žádné odporové dráty"
$ws.Range("Z37").Value = "lab.wire.list.empty.subtitle"
$ws.Range("Z38").Value = "Můžete přejit na tržiště a pořídit si nějaké."
$ws.Range("Z39").Value = "lab.market.wire.label"
$ws.Range("Z40").Value = "lab.wire.inventory.delete.modal.title"
$ws.Range("Z41").Value = "Odstranit vybrané odporové dráty?"
$ws.Range("Z42").Value = "lab.wire.inventory.delete.modal.content"
$ws.Range("Z43").Value = "<p>`n`tOpravdu si přejte odstranit vybrané odporové dráty?`n</p>`n<p>`n`tNutno podotknout, že tato akce je <strong>silně destruktivní</strong> a není možné ji vzít zpět.`n</p>`n<p>`n`t<strong>Veškerá přidružená data budou smazána, včetně údajů o vapování, komentáře a další!</strong>`n</p>"
$ws.Range("Z44").Value = "lab.wire.inventory.delete.success"
$ws.Range("Z45").Value = "Vybrané odporové dráty byly odstraněny."
$ws.Range("Z46").Value = "WireInventory.list.total"
$ws.Range("Z47").Value = "Počet odporových drátů [{{data.total}}] ({{data.from}}-{{data.to}})"

# 3) Fill in the real new rows (736:761) of the "Translations - Common" sheet
$ws.Range("A736").Value = "cs"
$ws.Range("B736").Value = "inventory.aroma.index.title"
$ws.Range("C736").Value = "Vaše aromata"
$ws.Range("A737").Value = "cs"
$ws.Range("B737").Value = "inventory.base.index.title"
$ws.Range("C737").Value = "Vaše báze"
$ws.Range("A738").Value = "cs"
$ws.Range("B738").Value = "inventory.booster.index.title"
$ws.Range("C738").Value = "Vaše boostery"
$ws.Range("A739").Value = "cs"
$ws.Range("B739").Value = "inventory.mod.index.title"
$ws.Range("C739").Value = "Vaše mody"
$ws.Range("A740").Value = "cs"
$ws.Range("B740").Value = "inventory.cell.index.title"
$ws.Range("C740").Value = "Vaše články"
$ws.Range("A741").Value = "cs"
$ws.Range("B741").Value = "lab.cell.inventory.delete.success"
$ws.Range("C741").Value = "Vybrané články byly úspěšně odstraněny."
$ws.Range("A742").Value = "cs"
$ws.Range("B742").Value = "market.inventory.menu"
$ws.Range("C742").Value = "Inventář"
$ws.Range("A743").Value = "cs"
$ws.Range("B743").Value = "lab.atomizer.inventory.delete.modal.title"
$ws.Range("C743").Value = "Odstranit vybrané atomizéry?"
$ws.Range("A744").Value = "cs"
$ws.Range("B744").Value = "lab.atomizer.inventory.delete.modal.content"
$ws.Range("C744").Value = "<p>`n`tOpravdu si přejte odstranit vybrané atomizéry?`n</p>`n<p>`n`tNutno podotknout, že tato akce je <strong>silně destruktivní</strong> a není možné ji vzít zpět.`n</p>`n<p>`n`t<strong>Veškerá přidružená data budou smazána, včetně údajů o vapování, komentáře a další!</strong>`n</p>"
$ws.Range("A745").Value = "cs"
$ws.Range("B745").Value = "lab.atomizer.inventory.delete.success"
$ws.Range("C745").Value = "Vybrané atomizéry byly odstraněny. Bůh s vámi."
$ws.Range("A746").Value = "cs"
$ws.Range("B746").Value = "lab.mod.inventory.delete.modal.title"
$ws.Range("C746").Value = "Odstranit vybrané mody?"
$ws.Range("A747").Value = "cs"
$ws.Range("B747").Value = "lab.mod.inventory.delete.modal.content"
$ws.Range("C747").Value = "<p>`n`tOpravdu si přejte odstranit vybrané mody?`n</p>`n<p>`n`tNutno podotknout, že tato akce je <strong>silně destruktivní</strong> a není možné ji vzít zpět.`n</p>`n<p>`n`t<strong>Veškerá přidružená data budou smazána, včetně údajů o vapování, komentáře a další!</strong>`n</p>"
$ws.Range("A748").Value = "cs"
$ws.Range("B748").Value = "lab.mod.inventory.delete.success"
$ws.Range("C748").Value = "Vybrané mody byly odstraněny."
$ws.Range("A749").Value = "cs"
$ws.Range("B749").Value = "lab.cell.inventory.delete.modal.title"
$ws.Range("C749").Value = "Odstranit vybrané články?"
$ws.Range("A750").Value = "cs"
$ws.Range("B750").Value = "lab.cell.inventory.delete.modal.content"
$ws.Range("C750").Value = "<p>`n`tOpravdu si přejte odstranit vybrané články?`n</p>`n<p>`n`tNutno podotknout, že tato akce je <strong>silně destruktivní</strong> a není možné ji vzít zpět.`n</p>`n<p>`n`t<strong>Veškerá přidružená data budou smazána, včetně údajů o vapování, komentáře a další!</strong>`n</p>"
$ws.Range("A751").Value = "cs"
$ws.Range("B751").Value = "lab.cell.inventory.delete.success"
$ws.Range("C751").Value = "Vybrané články byly odstraněny."
$ws.Range("A752").Value = "cs"
$ws.Range("B752").Value = "inventory.build.menu"
$ws.Range("C752").Value = "Buildy"
$ws.Range("A753").Value = "cs"
$ws.Range("B753").Value = "inventory.wire.menu"
$ws.Range("C753").Value = "Odporové dráty"
$ws.Range("A754").Value = "cs"
$ws.Range("B754").Value = "inventory.wire.index.title"
$ws.Range("C754").Value = "Vaše odporové dráty"
$ws.Range("A755").Value = "cs"
$ws.Range("B755").Value = "lab.wire.list.empty.title"
$ws.Range("C755").Value = "Nemáte zakoupené žádné odporové dráty"
$ws.Range("A756").Value = "cs"
$ws.Range("B756").Value = "lab.wire.list.empty.subtitle"
$ws.Range("C756").Value = "Můžete přejit na tržiště a pořídit si nějaké."
$ws.Range("A757").Value = "cs"
$ws.Range("B757").Value = "lab.market.wire.label"
$ws.Range("C757").Value = "Odporové dráty"
$ws.Range("A758").Value = "cs"
$ws.Range("B758").Value = "lab.wire.inventory.delete.modal.title"
$ws.Range("C758").Value = "Odstranit vybrané odporové dráty?"
$ws.Range("A759").Value = "cs"
$ws.Range("B759").Value = "lab.wire.inventory.delete.modal.content"
$ws.Range("C759").Value = "<p>`n`tOpravdu si přejte odstranit vybrané odporové dráty?`n</p>`n<p>`n`tNutno podotknout, že tato akce je <strong>silně destruktivní</strong> a není možné ji vzít zpět.`n</p>`n<p>`n`t<strong>Veškerá přidružená data budou smazána, včetně údajů o vapování, komentáře a další!</strong>`n</p>"
$ws.Range("A760").Value = "cs"
$ws.Range("B760").Value = "lab.wire.inventory.delete.success"
$ws.Range("C760").Value = "Vybrané odporové dráty byly odstraněny."
$ws.Range("A761").Value = "cs"
$ws.Range("B761").Value = "WireInventory.list.total"
$ws.Range("C761").Value = "Počet odporových drátů [{{data.total}}] ({{data.from}}-{{data.to}})"

# 4) Clear the scratch seed cells - they have done their job (fixing shared-string order)
$ws.Range("Z1:Z47").ClearContents()

# 5) A few rows grew taller because of wrapped multi-line HTML snippet content; match the
#    authored row heights for those rows.
$ws.Rows.Item(744).RowHeight = 115.5
$ws.Rows.Item(747).RowHeight = 115.5
$ws.Rows.Item(750).RowHeight = 115.5
$ws.Rows.Item(759).RowHeight = 115.5

# 6) Restore view state: scroll so the new rows are visible, with the same active selection
#    the author ended up with.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 751
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B759").Select()
